$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the label in F5 from "Sampling Period / Frame rate of camera (in seconds)"
# to "Sampling Rate / Frame rate of camera (in Hz)"
$ws.Range("F5").Value = "Sampling Rate / Frame rate of camera (in Hz)"

# Update the corresponding value in E5 from 1/60 (seconds) to 60 (Hz)
$ws.Range("E5").Value = 60

# Update selection to match final state (H5)
$ws.Range("H5").Select()
